$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, copying the header style (bold, border, centered) from E1
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       AdaBoostRegressor())]),`n                                            param_grid={'model__learning_rate': [0.1,`n                                                                                 0.5,`n                                                                                 1.0],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# Row 2: update metric values and add model description
$ws.Range("B2").Value = 0.1068557662344095
$ws.Range("C2").Value = 0.9980479982073243
$ws.Range("D2").Value = 0.2407779237252792
$ws.Range("F2").Value = $modelText

# Row 3: update metric values and add model description
$ws.Range("B3").Value = 0.3162511885970278
$ws.Range("C3").Value = 0.9766342100919515
$ws.Range("D3").Value = 0.4312949900126981
$ws.Range("F3").Value = $modelText
